# The "Stre" column (G) holding a stroke-width numeric value was set up
# incorrectly (commit: "fix: link color setting"); remove it entirely so
# that the following "lwd" column (H) shifts left into its place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(7).Delete()
